$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KeyValuePairs")

# Add the two new localization rows for the service-worker update modal.
# Written in the same order Excel would generate new shared-string entries:
# both keys first, then the EN/DE pair for each row.
$ws.Range("A48").Value = "service_worker-update_headline"
$ws.Range("A49").Value = "service_worker-update_confirm_btn_txt"

$ws.Range("B48").Value = "Update Available"
$ws.Range("C48").Value = "Update verfügbar"

$ws.Range("B49").Value = "Update Now & Refresh"
$ws.Range("C49").Value = "Update installieren"

# Match the formatting already used by the rest of the table (fill style
# carried by every non-wrapped data row) by copying it onto the new rows.
$ws.Range("A46:C46").Copy() | Out-Null
$ws.Range("A48:C49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Keep the table bound to the newly added rows.
$table = $ws.ListObjects.Item("Tabelle2")
$table.Resize($ws.Range("A1:C49"))
